$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I (I0) and J (IF), styled like the other
# header cells (bold / bordered / centered) by copying the style of H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Values for the new I0 / J0 (I/J) columns, rows 2..58
$iVals = @(6,5,7,9,5,8,11,8,7,8,7,4,8,7,6,5,7,7,8,7,7,9,8,6,7,7,7,6,7,8,9,6,7,8,6,7,8,8,7,8,8,8,7,7,8,9,8,8,8,7,6,6,6,5,5,9,5)
$jVals = @(6,5,7,9,5,8,12,8,7,8,7,4,8,7,7,6,7,8,9,7,7,9,9,6,7,7,8,7,7,8,9,6,7,8,6,7,8,8,7,8,8,8,8,7,8,9,8,8,8,7,6,6,6,5,5,9,5)

for ($r = 2; $r -le 58; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
